$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "28.373.15"
$ws.Range("E2").Value = "  +5.99%  "
Set-TextValue $ws.Range("D3") "1.813.83"
$ws.Range("E3").Value = "  +5.55%  "
$ws.Range("E4").Value = "  -0.24%  "
Set-TextValue $ws.Range("D5") "317.71"
$ws.Range("E5").Value = "  +2.83%  "
Set-TextValue $ws.Range("D6") "1.000"
$ws.Range("E6").Value = "  -0.22%  "
Set-TextValue $ws.Range("D7") "0.5738"
$ws.Range("E7").Value = "  +18.09%  "
Set-TextValue $ws.Range("D8") "0.3877"
$ws.Range("E8").Value = "  +11.57%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D9") "43.17"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D10") "0.07616"
$ws.Range("E10").Value = "  +5.67%  "
Set-TextValue $ws.Range("D11") "1.139"
$ws.Range("E11").Value = "  +9.15%  "
Set-TextValue $ws.Range("D12") "21.25"
$ws.Range("E12").Value = "  +7.79%  "
Set-TextValue $ws.Range("D13") "1.000"
$ws.Range("E13").Value = "  -0.27%  "
Set-TextValue $ws.Range("D14") "6.259"
$ws.Range("E14").Value = "  +7.14%  "
Set-TextValue $ws.Range("D15") "1.809.86"
$ws.Range("E15").Value = "  +4.85%  "
Set-TextValue $ws.Range("D16") "7.276"
$ws.Range("E16").Value = "  +7.18%  "
Set-TextValue $ws.Range("D17") "92.09"
$ws.Range("E17").Value = "  +7.00%  "
Set-TextValue $ws.Range("D18") "0.00001077"
$ws.Range("E18").Value = "  +4.39%  "
Set-TextValue $ws.Range("D19") "0.06478"
$ws.Range("E19").Value = "  +1.21%  "
Set-TextValue $ws.Range("D20") "0.9994"
$ws.Range("E20").Value = "  -0.27%  "
Set-TextValue $ws.Range("D21") "17.34"
$ws.Range("E21").Value = "  +5.27%  "
Set-TextValue $ws.Range("D22") "6.000"
$ws.Range("E22").Value = "  +5.46%  "
Set-TextValue $ws.Range("D23") "28.385.03"
$ws.Range("E23").Value = "  +5.74%  "
$ws.Range("E24").Value = "  +4.02%  "
Set-TextValue $ws.Range("D25") "2.122"
$ws.Range("E25").Value = "  +3.45%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D26") "158.41"
$ws.Range("E26").Value = "  +2.65%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D27") "20.85"
$ws.Range("E27").Value = "  +5.49%  "
Set-TextValue $ws.Range("D28") "2.438"
$ws.Range("E28").Value = "  +18.88%  "
Set-TextValue $ws.Range("D29") "2.020.55"
$ws.Range("E29").Value = "  +5.18%  "
Set-TextValue $ws.Range("D30") "124.03"
$ws.Range("E30").Value = "  +3.52%  "
Set-TextValue $ws.Range("D31") "1.167"
$ws.Range("E31").Value = "  +13.85%  "
Set-TextValue $ws.Range("D32") "0.1058"
$ws.Range("E32").Value = "  +14.03%  "
Set-TextValue $ws.Range("D33") "5.787"
$ws.Range("E33").Value = "  +8.47%  "
Set-TextValue $ws.Range("D34") "3.636"
$ws.Range("E34").Value = "  +2.11%  "
Set-TextValue $ws.Range("D35") "8.959"
$ws.Range("E35").Value = "  +21.03%  "
Set-TextValue $ws.Range("D36") "0.02323"
$ws.Range("E36").Value = "  +7.01%  "
Set-TextValue $ws.Range("D37") "0.2167"
$ws.Range("E37").Value = "  +9.63%  "
Set-TextValue $ws.Range("D38") "11.69"
$ws.Range("E38").Value = "  +7.46%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D39") "0.6436"
$ws.Range("E39").Value = "  +8.25%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D40") "0.06098"
$ws.Range("E40").Value = "  +3.90%  "
Set-TextValue $ws.Range("D41") "5.042"
$ws.Range("E41").Value = "  +6.99%  "
Set-TextValue $ws.Range("D42") "1.164"
$ws.Range("E42").Value = "  +4.44%  "
Set-TextValue $ws.Range("D43") "0.9996"
$ws.Range("E43").Value = "  -0.22%  "
Set-TextValue $ws.Range("D44") "1.375"
$ws.Range("E44").Value = "  -3.56%  "
$ws.Range("E45").Value = "  +4.93%  "
Set-TextValue $ws.Range("D46") "0.6006"
$ws.Range("E46").Value = "  +7.87%  "
Set-TextValue $ws.Range("D47") "3.707"
$ws.Range("E47").Value = "  +3.88%  "
Set-TextValue $ws.Range("D48") "122.90"
$ws.Range("E48").Value = "  +3.54%  "
$ws.Range("E49").Value = "  +6.74%  "
$ws.Range("E50").Value = "  +5.67%  "
Set-TextValue $ws.Range("D51") "0.06853"
$ws.Range("E51").Value = "  +3.59%  "
